# Generating GCODE from unrouted PCB.pptx - documentation changes and path fixes
#
# 1. Refresh the "datetimeFigureOut" auto date field cached on the slide
#    master and on every slide layout (6/28/2023 -> 8/10/2023).
# 2. Slide 2 ("Requirements"): reword the KiCAD netlist requirement and
#    fix the capitalisation of "EEschema".

$p = $ppt.ActivePresentation

# --- 1. Update the cached "Date Placeholder" text everywhere it lives ---

function Update-DatePlaceholder($container, [string]$newDate) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $len = $tr.Length
            if ($len -gt 0) {
                $tr.Characters(1, $len).Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master "8/10/2023"

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($L) "8/10/2023"
}

# --- 2. Slide 2 content tweaks -----------------------------------------

$slide = $p.Slides.Item(2)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$full = $tr.Text
$old = "Valid netlist from "
$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $tr.Characters($idx + 1, $old.Length).Text = "Valid netlist compatible with "
}

$full = $tr.Text
$old = "Eeschema"
$idx = $full.IndexOf($old)
if ($idx -ge 0) {
    $tr.Characters($idx + 1, $old.Length).Text = "EEschema"
}

Write-Output "done"
